$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 8239.704
$ws.Range("I98").Value = 8717.378000000001
$ws.Range("J98").Value = 5714.857
$ws.Range("K98").Value = 8717.378000000001
$ws.Range("L98").Value = 5714.857
$ws.Range("M98").Value = -7219.378000000001
$ws.Range("N98").Value = -8710.857

$ws.Range("H106").Value = 3052
$ws.Range("I106").Value = 2819.4375
$ws.Range("K106").Value = 2819.4375
$ws.Range("M106").Value = -2188.4375

$ws.Range("H122").Value = 8239.704
$ws.Range("I122").Value = 8717.378000000001
$ws.Range("J122").Value = 5714.857
$ws.Range("K122").Value = 26152.134
$ws.Range("L122").Value = 17144.571
$ws.Range("M122").Value = -23702.134
$ws.Range("N122").Value = -22044.571

$ws.Range("H138").Value = 1312512.2
$ws.Range("I138").Value = 11193.5
$ws.Range("J138").Value = 1526427.6
$ws.Range("K138").Value = 33580.5
$ws.Range("L138").Value = 4579282.800000001
$ws.Range("M138").Value = -28440.5
$ws.Range("N138").Value = -4589562.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 727.3200000000001
$ws.Range("I97").Value = 495.6842
$ws.Range("K97").Value = 495.6842
$ws.Range("M97").Value = 0.3158000000000243

$ws.Range("H133").Value = 74978.3
$ws.Range("J133").Value = 74978.3
$ws.Range("L133").Value = 74978.3
$ws.Range("N133").Value = -80038.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1250
$ws.Range("I29").Value = 1250
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1250
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -961
$ws.Range("N29").Value = ""

$ws.Range("H82").Value = 71468300
$ws.Range("I82").Value = 34583.668
$ws.Range("K82").Value = 34583.668
$ws.Range("M82").Value = -34200.668

$ws.Range("H85").Value = 71468300
$ws.Range("I85").Value = 34583.668
$ws.Range("K85").Value = 34583.668
$ws.Range("M85").Value = -33257.668

$ws.Range("H94").Value = 2232.1538
$ws.Range("I94").Value = 2232.1538
$ws.Range("K94").Value = 2232.1538
$ws.Range("M94").Value = -1781.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 19585.572
$ws.Range("I69").Value = 14420.8
$ws.Range("J69").Value = 32497.5
$ws.Range("K69").Value = 14420.8
$ws.Range("L69").Value = 32497.5
$ws.Range("M69").Value = -13671.8
$ws.Range("N69").Value = -33995.5

$ws.Range("H72").Value = 19585.572
$ws.Range("I72").Value = 14420.8
$ws.Range("J72").Value = 32497.5
$ws.Range("K72").Value = 43262.39999999999
$ws.Range("L72").Value = 97492.5
$ws.Range("M72").Value = -39518.39999999999
$ws.Range("N72").Value = -104980.5

$ws.Range("H86").Value = 3013.0454
$ws.Range("I86").Value = 2235.6428
$ws.Range("K86").Value = 2235.6428
$ws.Range("M86").Value = -1112.6428

$ws.Range("H89").Value = 3013.0454
$ws.Range("I89").Value = 2235.6428
$ws.Range("K89").Value = 11178.214
$ws.Range("M89").Value = -5562.214

$ws.Range("H93").Value = 12785.111
$ws.Range("I93").Value = 8133.25
$ws.Range("J93").Value = 50000
$ws.Range("K93").Value = 8133.25
$ws.Range("L93").Value = 50000
$ws.Range("M93").Value = -6261.25
$ws.Range("N93").Value = -53744

$ws.Range("H99").Value = 3375
$ws.Range("I99").Value = 2750
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2750
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -1252
$ws.Range("N99").Value = -6996

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = ""

$ws.Range("H122").Value = 3436.0952
$ws.Range("I122").Value = 3432.9
$ws.Range("K122").Value = 10298.7
$ws.Range("M122").Value = -7848.700000000001

$ws.Range("H126").Value = 3375
$ws.Range("I126").Value = 2750
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8250
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -5780
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 1585.7241
$ws.Range("I132").Value = 1359.48
$ws.Range("J132").Value = 2999.75
$ws.Range("K132").Value = 4078.44
$ws.Range("L132").Value = 8999.25
$ws.Range("M132").Value = -1548.44
$ws.Range("N132").Value = -14059.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 228
$ws.Range("I33").Value = 281.625
$ws.Range("J33").Value = 180.33333
$ws.Range("K33").Value = 1689.75
$ws.Range("L33").Value = 1081.99998
$ws.Range("M33").Value = -1406.75
$ws.Range("N33").Value = -1647.99998

$ws.Range("H40").Value = 208.05882
$ws.Range("J40").Value = 1025.6666
$ws.Range("L40").Value = 4102.6664
$ws.Range("N40").Value = -4240.6664

$ws.Range("H132").Value = 3129
$ws.Range("J132").Value = 3124.6875
$ws.Range("L132").Value = 28122.1875
$ws.Range("N132").Value = -33182.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 250038130
$ws.Range("J34").Value = 250038130
$ws.Range("L34").Value = 250038130
$ws.Range("N34").Value = -250038666

$ws.Range("H76").Value = 250038130
$ws.Range("J76").Value = 250038130
$ws.Range("L76").Value = 250038130
$ws.Range("N76").Value = -250038760

$ws.Range("H79").Value = 250038130
$ws.Range("J79").Value = 250038130
$ws.Range("L79").Value = 250038130
$ws.Range("N79").Value = -250040314

$ws.Range("H80").Value = 7912.25
$ws.Range("J80").Value = 8756.857
$ws.Range("L80").Value = 8756.857
$ws.Range("N80").Value = -10752.857

$ws.Range("H83").Value = 7912.25
$ws.Range("J83").Value = 8756.857
$ws.Range("L83").Value = 43784.285
$ws.Range("N83").Value = -53768.285

$ws.Range("H97").Value = 1183.5834
$ws.Range("I97").Value = 934.8889
$ws.Range("J97").Value = 1929.6666
$ws.Range("K97").Value = 934.8889
$ws.Range("L97").Value = 1929.6666
$ws.Range("M97").Value = -438.8889
$ws.Range("N97").Value = -2921.6666

$ws.Range("H102").Value = 31252804
$ws.Range("J102").Value = 5854.4
$ws.Range("L102").Value = 5854.4
$ws.Range("N102").Value = -9098.4

$ws.Range("H122").Value = 4201.636
$ws.Range("I122").Value = 1459.7142
$ws.Range("K122").Value = 4379.142599999999
$ws.Range("M122").Value = -1929.142599999999

$ws.Range("H132").Value = 2923
$ws.Range("I132").Value = 2559.4375
$ws.Range("J132").Value = 3130.75
$ws.Range("K132").Value = 7678.3125
$ws.Range("L132").Value = 9392.25
$ws.Range("M132").Value = -5148.3125
$ws.Range("N132").Value = -14452.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2844.5557
$ws.Range("I7").Value = 2950.75
$ws.Range("K7").Value = 2950.75
$ws.Range("M7").Value = -2838.75

$ws.Range("H22").Value = 2609.4
$ws.Range("I22").Value = 865.6667
$ws.Range("K22").Value = 865.6667
$ws.Range("M22").Value = -570.6667

$ws.Range("H27").Value = 2609.4
$ws.Range("I27").Value = 865.6667
$ws.Range("K27").Value = 865.6667
$ws.Range("M27").Value = -758.6667

$ws.Range("H40").Value = 5220.273
$ws.Range("I40").Value = 3740.25
$ws.Range("K40").Value = 3740.25
$ws.Range("M40").Value = -3604.25

$ws.Range("H55").Value = 1319.4
$ws.Range("J55").Value = 1524.5
$ws.Range("L55").Value = 1524.5
$ws.Range("N55").Value = -1870.5

$ws.Range("H126").Value = 2844.5557
$ws.Range("I126").Value = 2950.75
$ws.Range("K126").Value = 8852.25
$ws.Range("M126").Value = -6382.25

$ws.Range("H132").Value = 1338112
$ws.Range("J132").Value = 3083878.8
$ws.Range("L132").Value = 9251636.399999999
$ws.Range("N132").Value = -9256696.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 4157.143
$ws.Range("I52").Value = 4157.143
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 4157.143
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -3931.143
$ws.Range("N52").Value = ""

$ws.Range("H81").Value = 3990.0908
$ws.Range("J81").Value = 6999.5
$ws.Range("L81").Value = 13999
$ws.Range("N81").Value = -16121

$ws.Range("H84").Value = 3990.0908
$ws.Range("J84").Value = 6999.5
$ws.Range("L84").Value = 69995
$ws.Range("N84").Value = -80603

$ws.Range("H98").Value = 113333.336
$ws.Range("J98").Value = 113333.336
$ws.Range("L98").Value = 113333.336
$ws.Range("N98").Value = -119323.336

$ws.Range("H122").Value = 3570.889
$ws.Range("I122").Value = 1666
$ws.Range("K122").Value = 4998
$ws.Range("M122").Value = -2548

$ws.Range("H126").Value = 2969.8
$ws.Range("I126").Value = 3087.25
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 9261.75
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -6791.75
$ws.Range("N126").Value = -12440

$ws.Range("H132").Value = 6480.9165
$ws.Range("J132").Value = 7506.5
$ws.Range("L132").Value = 22519.5
$ws.Range("N132").Value = -27579.5
